# Excel COM-interop edit script
# Applies the scheduled-runner numeric corrections to the H:N price/profit
# columns across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H17").Value = 3486.7834
$ws.Range("J17").Value = 3486.7834
$ws.Range("L17").Value = 10460.3502
$ws.Range("N17").Value = -10796.3502

$ws.Range("H33").Value = 302.66666
$ws.Range("I33").Value = 344.75
$ws.Range("J33").Value = 218.5
$ws.Range("K33").Value = 344.75
$ws.Range("L33").Value = 218.5
$ws.Range("M33").Value = -115.75
$ws.Range("N33").Value = -676.5

$ws.Range("H42").Value = 155.57143
$ws.Range("I42").Value = 155.57143
$ws.Range("K42").Value = 466.71429
$ws.Range("M42").Value = -236.71429

$ws.Range("H51").Value = 21041.479
$ws.Range("I51").Value = 11795.7
$ws.Range("J51").Value = 28153.615
$ws.Range("K51").Value = 11795.7
$ws.Range("L51").Value = 28153.615
$ws.Range("M51").Value = -11311.7
$ws.Range("N51").Value = -29121.615

$ws.Range("H70").Value = 1372665.6
$ws.Range("J70").Value = 2118664.8
$ws.Range("L70").Value = 6355994.399999999
$ws.Range("N70").Value = -6356534.399999999

$ws.Range("H73").Value = 1372665.6
$ws.Range("J73").Value = 2118664.8
$ws.Range("L73").Value = 6355994.399999999
$ws.Range("N73").Value = -6357866.399999999

$ws.Range("H76").Value = 9093613
$ws.Range("I76").Value = 14287764
$ws.Range("K76").Value = 14287764
$ws.Range("M76").Value = -14287449

$ws.Range("H79").Value = 9093613
$ws.Range("I79").Value = 14287764
$ws.Range("K79").Value = 14287764
$ws.Range("M79").Value = -14286672

$ws.Range("H132").Value = 2107.3572
$ws.Range("I132").Value = 1275.3334
$ws.Range("K132").Value = 3826.0002
$ws.Range("M132").Value = -1296.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1938.9
$ws.Range("I45").Value = 1823.875
$ws.Range("K45").Value = 1823.875
$ws.Range("M45").Value = -1446.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 88402.17
$ws.Range("I86").Value = 1465.45
$ws.Range("K86").Value = 1465.45
$ws.Range("M86").Value = -342.45

$ws.Range("H89").Value = 88402.17
$ws.Range("I89").Value = 1465.45
$ws.Range("K89").Value = 7327.25
$ws.Range("M89").Value = -1711.25

$ws.Range("H94").Value = 1560.8667
$ws.Range("J94").Value = 1432.9333
$ws.Range("L94").Value = 1432.9333
$ws.Range("N94").Value = -2334.9333

$ws.Range("H107").Value = 1649.7916
$ws.Range("I107").Value = 1318
$ws.Range("K107").Value = 1318
$ws.Range("M107").Value = 602

$ws.Range("H132").Value = 110684
$ws.Range("J132").Value = 110684
$ws.Range("L132").Value = 110684
$ws.Range("N132").Value = -120804

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3046.1
$ws.Range("I4").Value = 103.333336
$ws.Range("J4").Value = 4307.2856
$ws.Range("K4").Value = 103.333336
$ws.Range("L4").Value = 4307.2856
$ws.Range("M4").Value = 8.666663999999997
$ws.Range("N4").Value = -4531.2856

$ws.Range("H29").Value = 10186.5
$ws.Range("I29").Value = 6719
$ws.Range("J29").Value = 10880
$ws.Range("K29").Value = 6719
$ws.Range("L29").Value = 10880
$ws.Range("M29").Value = -6426
$ws.Range("N29").Value = -11466

$ws.Range("H58").Value = 8259
$ws.Range("I58").Value = 5012
$ws.Range("K58").Value = 5012
$ws.Range("M58").Value = -4809

$ws.Range("H136").Value = 8259
$ws.Range("I136").Value = 5012
$ws.Range("K136").Value = 15036
$ws.Range("M136").Value = -12486

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8500
$ws.Range("I3").Value = 8500
$ws.Range("K3").Value = 25500
$ws.Range("M3").Value = -25388

$ws.Range("H37").Value = 95286.586
$ws.Range("J37").Value = 95286.586
$ws.Range("L37").Value = 285859.758
$ws.Range("N37").Value = -286083.758

$ws.Range("H81").Value = 2464.3333
$ws.Range("I81").Value = 2464.3333
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 7392.999899999999
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -6269.999899999999
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 2464.3333
$ws.Range("I84").Value = 2464.3333
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 22178.9997
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -16562.9997
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7700.3125
$ws.Range("I132").Value = 3843.7144
$ws.Range("K132").Value = 11531.1432
$ws.Range("M132").Value = -9001.143199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4403.5947
$ws.Range("I22").Value = 2479.2856
$ws.Range("J22").Value = 5574.913
$ws.Range("K22").Value = 2479.2856
$ws.Range("L22").Value = 5574.913
$ws.Range("M22").Value = -2184.2856
$ws.Range("N22").Value = -6164.913

$ws.Range("H27").Value = 4403.5947
$ws.Range("I27").Value = 2479.2856
$ws.Range("J27").Value = 5574.913
$ws.Range("K27").Value = 2479.2856
$ws.Range("L27").Value = 5574.913
$ws.Range("M27").Value = -2372.2856
$ws.Range("N27").Value = -5788.913

$ws.Range("H40").Value = 20642.285
$ws.Range("I40").Value = 17299.4
$ws.Range("K40").Value = 17299.4
$ws.Range("M40").Value = -17163.4

$ws.Range("H55").Value = 138.95653
$ws.Range("J55").Value = 159.83333
$ws.Range("L55").Value = 159.83333
$ws.Range("N55").Value = -505.83333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4762.8335
$ws.Range("I81").Value = 4377.6816
$ws.Range("K81").Value = 8755.3632
$ws.Range("M81").Value = -7694.3632

$ws.Range("H84").Value = 4762.8335
$ws.Range("I84").Value = 4377.6816
$ws.Range("K84").Value = 43776.816
$ws.Range("M84").Value = -38472.816

$ws.Range("H94").Value = 45697.5
$ws.Range("J94").Value = 45697.5
$ws.Range("L94").Value = 45697.5
$ws.Range("N94").Value = -47499.5

$ws.Range("H113").Value = 1465.25
$ws.Range("I113").Value = 1531.7142
$ws.Range("K113").Value = 4595.142599999999
$ws.Range("M113").Value = -2425.142599999999

$ws.Range("H122").Value = 3668.25
$ws.Range("I122").Value = 3642.647
$ws.Range("K122").Value = 10927.941
$ws.Range("M122").Value = -8477.940999999999

$ws.Range("H136").Value = 6044.1665
$ws.Range("I136").Value = 2738.5
$ws.Range("K136").Value = 8215.5
$ws.Range("M136").Value = -5665.5
